$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 83.5
$ws.Range("C3").Value = 83.59999999999999
$ws.Range("C4").Value = 82.2
$ws.Range("C5").Value = 89.09999999999999
$ws.Range("C6").Value = 89.2
$ws.Range("C7").Value = 89.5
